$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "93.317.64"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("D3").Value = "3.128.28"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.50"
$ws.Range("E5").Value = "  -2.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "614.94"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  +2.52%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.815"
$ws.Range("E10").Value = "  +8.29%  "
$ws.Range("D11").Value = "3.125.39"
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("E12").Value = "  -1.94%  "
$ws.Range("E13").Value = "  -2.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.90"
$ws.Range("E14").Value = "  +0.16%  "
$ws.Range("D15").Value = "93.025.04"
$ws.Range("E15").Value = "  +1.66%  "
$ws.Range("E16").Value = "  -2.82%  "
$ws.Range("D17").Value = "3.711.30"
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("D18").Value = "3.135.03"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.80"
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.67"
$ws.Range("E20").Value = "  -1.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.96"
$ws.Range("E21").Value = "  +2.83%  "
$ws.Range("E22").Value = "  +0.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "442.99"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.17"
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.22"
$ws.Range("E25").Value = "  +5.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.65"
$ws.Range("E26").Value = "  -3.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.50"
$ws.Range("E27").Value = "  +6.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "85.98"
$ws.Range("E28").Value = "  -3.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("E30").Value = "  +7.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.233"
$ws.Range("E31").Value = "  +3.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.123"
$ws.Range("E32").Value = "  -10.22%  "
$ws.Range("E33").Value = "  -1.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.03"
$ws.Range("E34").Value = "  -1.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.15"
$ws.Range("E35").Value = "  +5.97%  "
$ws.Range("E36").Value = "  -9.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "25.92"
$ws.Range("E37").Value = "  -1.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.96"
$ws.Range("E38").Value = "  +0.73%  "
$ws.Range("E39").Value = "  -2.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.30"
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.99"
$ws.Range("E41").Value = "  +7.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "475.70"
$ws.Range("E42").Value = "  -2.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.433"
$ws.Range("E43").Value = "  -0.55%  "
$ws.Range("E44").Value = "  -2.63%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "159.17"
$ws.Range("E46").Value = "  -0.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.692"
$ws.Range("E47").Value = "  -1.00%  "
$ws.Range("E48").Value = "  -3.07%  "
$ws.Range("E49").Value = "  -1.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "44.11"
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.41"
$ws.Range("E51").Value = "  -0.53%  "
